$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values to apply to rows 2 and 3 (same values for both rows)
$updates = @{
    "D" = -0.0117
    "E" = 0.0383
    "F" = 0.0275
    "G" = 0.06510033964864353
    "H" = 0.06510033964864353
    "I" = 0.05826671463619382
    "J" = 0.04065616292975397
    "K" = 1576.7
    "L" = 0.03204811587486458
    "W" = 0.07626413597624092
    "X" = 0.1332976006534351
    "Y" = -0.05703346467719418
    "Z" = 1.162221922992424
    "AA" = 0.04725148386171196
    "AB" = 0.05611052168252498
    "AC" = -0.008859037820813022
    "AD" = 22602.2
    "AF" = 22602.2
    "AG" = 22602.2
    "AH" = 0.6714714786085851
    "AI" = 0.4848135151135554
    "AJ" = 0.6714714786085851
    "AK" = 0.4848135151135554
    "AL" = 315.4
    "AM" = 315.4
    "AN" = 7.614783370392831
    "AO" = 9.088776157260622
    "AP" = 7.614783370392831
    "AQ" = 9.088776157260622
}

foreach ($row in 2, 3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
